# "fixed horizontal centering on registers"
#
# The register sheet (СТР) is printed with the table not centered on the
# page. Turn on horizontal centering for printing, and leave the
# selection where the edit was made (the registration-number column of
# the header table, C9:C11) as the last interactively-selected range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center the printed sheet horizontally on the page (adds
# <printOptions horizontalCentered="1"/> to the worksheet).
$ws.PageSetup.CenterHorizontally = $true

# Leave the selection on the merged "№ п/п" column of the roster header
# (C9:C11), matching where the centering fix was authored.
$ws.Range("C9:C11").Select()
